# Daily attendance processing - 2025-10-31 15:46:03
# Normalize the "Recorded By" (column G) entries: swap the order of the
# first two comma-separated recorders (e.g. "System, X" -> "X, System")
# for every row whose recorder list begins with "System, " and whose
# second recorder is not "backup@backdoor.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (in column G) that need the first two comma-separated tokens swapped.
$targetRows = @(2,3,6,7,10,11,12,13,14,15,17,18,19,20,21,22,24,29,30,33,34,37,38,39,40,41,42,44,45,46,47,48,49,51,56,57,60,61,64,65,66,67,68,69,71,72,73,74,75,76,78,86,87,88,89,93,95,96,97,99,102,104,112,113,114,115,119,121,122,123,125,128,130,138,139,140,141,145,147,148,149,151,154,156)

foreach ($row in $targetRows) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Text

    if ($val -eq $null) { continue }
    if (-not $val.StartsWith("System, ")) { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }
    if ($parts[1] -eq "backup@backdoor.com") { continue }

    $tmp = $parts[0]
    $parts[0] = $parts[1]
    $parts[1] = $tmp

    $cell.Value = ($parts -join ", ")
}
